## BOT; UPDATE DATA
## Adds two new daily rows (2020-04-14 / 2020-04-15) to the "相談件数"
## sheet just above the trailing note row, then fixes up the print area
## and the saved view state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")
$ws.Activate()

# Insert two fresh rows right before the old last row (the note/total
# row at row 80) -- this pushes that row down to 82 and inherits the
# number formats/styles from the row above, exactly like Excel's own
# "Insert Copied/Above" row behaviour.
$ws.Rows("80:81").Insert()

# New daily figures.
$ws.Range("A80").Value = 43935
$ws.Range("B80").Value = 770
$ws.Range("C80").Value = 22775
$ws.Range("D80").Value = 165
$ws.Range("E80").Value = 5232

$ws.Range("A81").Value = 43936
$ws.Range("B81").Value = 717
$ws.Range("C81").Value = 23492
$ws.Range("D81").Value = 138
$ws.Range("E81").Value = 5370

# The print area grows by the same two rows (was $A$1:$E$84).
$n = $wb.Names.Item("相談件数!Print_Area")
$n.RefersTo = '=相談件数!$A$1:$E$86'

# Keep the saved selection in sync with the new bottom of the sheet.
$ws.Range("E82").Select()
